# Marksheet update: fill in the student's answers (col A) for the quiz,
# fold the three side-by-side "Student Ans/Correct Ans" blocks into one,
# and refresh the score summary (rows 10-12) to reflect the graded totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- helper: Union of a list of A1 cell refs on $ws --------------------
function Get-UnionRange([string[]]$refs) {
    $u = $ws.Range($refs[0])
    for ($i = 1; $i -lt $refs.Length; $i++) {
        $u = $excel.Union($u, $ws.Range($refs[$i]))
    }
    return $u
}

# =========================================================================
# 1) Score summary block (rows 10-12): give the header cells in column A
#    the same "mtitleStyle" look as the row above, and write the graded
#    counts / marking scheme / totals.
# =========================================================================
$ws.Range("A9").Copy() | Out-Null
(Get-UnionRange @("A10", "A11", "A12")).PasteSpecial($xlPasteFormats) | Out-Null

# Right / Wrong / Not-attempted / Max counts
$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 28

# Marking scheme (points per right / wrong answer)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Totals
$ws.Range("B12").Value = 72
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "69/112"

# =========================================================================
# 2) Student answers (column A, rows 16-40) — fill in what the student
#    actually picked, then colour-code it against the correct answer
#    (column B) using the workbook's existing named styles:
#      correctStyle   (green) -> answer matches the correct one
#      incorrectStyle (red)   -> answer given but wrong
#      normalStyle    (black) -> left blank / not attempted
# =========================================================================
$answers = @{
    16 = "Option A"; 17 = "Option D"; 18 = "Option B"; 19 = "Option C";
    20 = "Option B"; 21 = "Option C"; 22 = "Option D"; 23 = "Option D";
    25 = "Option A"; 28 = "Option D"; 30 = "Option B"; 31 = "Option D";
    32 = "Option C"; 33 = "Option A"; 34 = "Option A"; 35 = "Option D";
    36 = "Option A"; 39 = "Option A"
}
foreach ($row in $answers.Keys) {
    $ws.Range("A$row").Value = $answers[$row]
}

$correctRows   = @(16,17,18,19,20,21,22,23,25,28,30,31,32,35,36)
$incorrectRows = @(33,34,39)
$blankRows     = @(24,26,27,29,37,38,40)

$ws.Range("B10").Copy() | Out-Null
(Get-UnionRange ($correctRows | ForEach-Object { "A$_" })).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("C10").Copy() | Out-Null
(Get-UnionRange ($incorrectRows | ForEach-Object { "A$_" })).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("D10").Copy() | Out-Null
(Get-UnionRange ($blankRows | ForEach-Object { "A$_" })).PasteSpecial($xlPasteFormats) | Out-Null

# =========================================================================
# 3) The second question block (cols D/E) is being retired now that col A
#    carries the answers directly — its first three rows (16-18) still
#    mirror column A for one more release, the rest (19-40) is cleared.
# =========================================================================
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"

$ws.Range("B10").Copy() | Out-Null
(Get-UnionRange @("D16","D17","D18")).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("D19:E40").Clear() | Out-Null

# =========================================================================
# 4) The third question block (cols G/H) is removed entirely.
# =========================================================================
$ws.Range("G15:H21").Clear() | Out-Null
